$d = $word.ActiveDocument

# Make sure edits are applied directly (not as tracked-change insert/delete markup)
$d.TrackRevisions = $false

# --- Edit 1 -------------------------------------------------------------
# "* Título del recurso (65 caracteres máx.) Refuerza tu aprendizaje: Los
#  Derechos del Hombre y del Ciudadano" -> drop the "Refuerza tu
#  aprendizaje: " prefix, keeping everything else (including the leading
#  space) intact in the same run.
$d.Content.Find.Execute(
    "Refuerza tu aprendizaje: Los Derechos del Hombre y del Ciudadano",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Los Derechos del Hombre y del Ciudadano", 2) | Out-Null

# --- Edit 2 ---------------------------------------------------------------
# Same visible text change on the second occurrence ("* Título del
# ejercicio ...") but this time Word's edit leaves the auto "_GoBack"
# bookmark right where the deletion happened, splitting the run in two:
#   <run>" "</run><bookmarkStart/><bookmarkEnd/><run>"Los Derechos..."</run>
$target = $d.Content
$found = $target.Find.Execute(
    "Refuerza tu aprendizaje: Los Derechos del Hombre y del Ciudadano",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $deleteRange = $d.Range($target.Start, $target.Start + 26)
    $deleteRange.Text = ""

    # Move/create the "_GoBack" bookmark to the collapsed point left behind
    # by the deletion -- this also removes it from its previous location
    # (end of the "Explicación (500 caracteres máximo)" paragraph) because
    # a document may only have one bookmark with a given name.
    $goBackRange = $d.Range($deleteRange.Start, $deleteRange.Start)
    $d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
}
